$wb = $excel.ActiveWorkbook

# Helper: write a literal TEXT value into a cell (preserving shared-string /
# text type even when the text looks like a plain number, e.g. "0.79").
# A direct `$range.Value = "0.79"` would be auto-coerced to a numeric cell
# by Excel, which is not what we want here (the workbook stores these as
# text). Instead we write a text formula that evaluates to the desired
# string, then convert that formula to a static value via copy / paste-
# special-values. This keeps the cell's declared type as text without
# introducing any new number formats / cell styles.
function Set-TextValue {
    param($range, [string]$text)
    $escaped = $text.Replace('"', '""')
    $range.Formula = '="' + $escaped + '"'
    $range.Copy()
    $range.PasteSpecial(-4163)  # xlPasteValues
}

# --- Restricciones_del_follower ---------------------------------------
$ws = $wb.Worksheets.Item("Restricciones_del_follower")

# Row 2 (J_0_L0_v)
Set-TextValue $ws.Range("A2") "-13.522157494894248 - 2x_1 + 5.618111148463803y_1 - 0.747506399099227y_2"
Set-TextValue $ws.Range("B2") "16.02215749489425"
Set-TextValue $ws.Range("D2") "0.79"
Set-TextValue $ws.Range("E2") "1.3"
Set-TextValue $ws.Range("F2") "9.8"

# Row 3 (J_0_LP_v)
Set-TextValue $ws.Range("A3") "-3.0343920415696015 + x_1 - 3x_2 - 0.09592779602473588y_1 + 0.9947551815306963y_2"
Set-TextValue $ws.Range("B3") "1.0343920415696015"
Set-TextValue $ws.Range("D3") "0.09"
Set-TextValue $ws.Range("E3") "2.9"
Set-TextValue $ws.Range("F3") "1.5"

# Row 4 (J_Ne_L0_v)
Set-TextValue $ws.Range("A4") "-35.66331209314874 + x_1 + x_2 + 4.936467928230815y_1 + 0.2698996457339049y_2"
Set-TextValue $ws.Range("B4") "33.28331209314874"
Set-TextValue $ws.Range("D4") "0.54"
Set-TextValue $ws.Range("E4") "2.3000000000000003"
Set-TextValue $ws.Range("F4") "8.7"

# --- Punto_modificado ----------------------------------------------------
$ws = $wb.Worksheets.Item("Punto_modificado")
Set-TextValue $ws.Range("A2") "6.65"
Set-TextValue $ws.Range("B2") "2.1"
Set-TextValue $ws.Range("C2") "5.2"
Set-TextValue $ws.Range("D2") "3.2"

# --- Vector_bf -------------------------------------------------------------
$ws = $wb.Worksheets.Item("Vector_bf")
Set-TextValue $ws.Range("A2") "-3.095366986888819"
Set-TextValue $ws.Range("A3") "-0.644743719745682"

# --- Vector_BF ---------------------------------------------------------
$ws = $wb.Worksheets.Item("Vector_BF")
Set-TextValue $ws.Range("A2") "-0.6000000000000001"
Set-TextValue $ws.Range("A3") "5.399999999999999"
Set-TextValue $ws.Range("A4") "-18.879230119462086"
Set-TextValue $ws.Range("A5") "-2.5338008927980056"

# --- Vector_Alpha --------------------------------------------------------
# These two cells are genuine numeric cells (no shared-string indirection),
# so a plain numeric assignment is correct here.
$ws = $wb.Worksheets.Item("Vector_Alpha")
$ws.Range("A2").Value = 1.757141052833688
$ws.Range("A3").Value = 0.09607106833453702
